# Update the Ptn-Alk NATMI TPM output sheet with the new TPM-based numbers,
# and add the "Resolving-Mac" sending cluster rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write one full data row (columns A:T) given row index and values
# ---------------------------------------------------------------------------
function Set-Row($r, $A, $B, $C, $D, $E, $F, $G, $H, $I, $J, $K, $L, $M, $N, $O, $P, $Q, $R, $S, $T) {
    $ws.Cells.Item($r, 1).Value  = $A
    $ws.Cells.Item($r, 2).Value  = $B
    $ws.Cells.Item($r, 3).Value  = $C
    $ws.Cells.Item($r, 4).Value  = $D
    $ws.Cells.Item($r, 5).Value  = $E
    $ws.Cells.Item($r, 6).Value  = $F
    $ws.Cells.Item($r, 7).Value  = $G
    $ws.Cells.Item($r, 8).Value  = $H
    $ws.Cells.Item($r, 9).Value  = $I
    $ws.Cells.Item($r, 10).Value = $J
    $ws.Cells.Item($r, 11).Value = $K
    $ws.Cells.Item($r, 12).Value = $L
    $ws.Cells.Item($r, 13).Value = $M
    $ws.Cells.Item($r, 14).Value = $N
    $ws.Cells.Item($r, 15).Value = $O
    $ws.Cells.Item($r, 16).Value = $P
    $ws.Cells.Item($r, 17).Value = $Q
    $ws.Cells.Item($r, 18).Value = $R
    $ws.Cells.Item($r, 19).Value = $S
    $ws.Cells.Item($r, 20).Value = $T
}

# Row 2: ECs -> Ptn -> Alk -> FAPs
Set-Row 2 "ECs" "Ptn" "Alk" "FAPs" `
    3 1 2.327816333333333 6.983449 `
    0.02128501190197005 0.02128501190197004 `
    1 0.3333333333333333 0.04487666666666667 0.13463 `
    0.7792485920506572 0.7792485920506572 `
    0.1044646376522222 0.94018173887 `
    0.01658631555639164 0.01658631555639164

# Row 3: ECs -> Ptn -> Alk -> MuSCs
Set-Row 3 "ECs" "Ptn" "Alk" "MuSCs" `
    3 1 2.327816333333333 6.983449 `
    0.02128501190197005 0.02128501190197004 `
    2 0.6666666666666666 0.012713 0.03813900000000001 `
    0.2207514079493428 0.2207514079493428 `
    0.02959352904566667 0.2663417614110001 `
    0.004698696345578407 0.004698696345578406

# Row 4: FAPs -> Ptn -> Alk -> FAPs
Set-Row 4 "FAPs" "Ptn" "Alk" "FAPs" `
    3 1 103.0385286666667 309.115586 `
    0.9421603747796319 0.9421603747796318 `
    1 0.3333333333333333 0.04487666666666667 0.13463 `
    0.7792485920506572 0.7792485920506572 `
    4.624025704797778 41.61623134318 `
    0.7341771455329477 0.7341771455329476

# Row 5: FAPs -> Ptn -> Alk -> MuSCs
Set-Row 5 "FAPs" "Ptn" "Alk" "MuSCs" `
    3 1 103.0385286666667 309.115586 `
    0.9421603747796319 0.9421603747796318 `
    2 0.6666666666666666 0.012713 0.03813900000000001 `
    0.2207514079493428 0.2207514079493428 `
    1.309928814939334 11.789359334454 `
    0.2079832292466842 0.2079832292466842

# Row 6: MuSCs -> Ptn -> Alk -> FAPs
Set-Row 6 "MuSCs" "Ptn" "Alk" "FAPs" `
    3 1 3.975769333333333 11.927308 `
    0.03635351138648862 0.03635351138648861 `
    1 0.3333333333333333 0.04487666666666667 0.13463 `
    0.7792485920506572 0.7792485920506572 `
    0.1784192751155556 1.60577347604 `
    0.02832842256401879 0.02832842256401879

# Row 7: MuSCs -> Ptn -> Alk -> MuSCs
Set-Row 7 "MuSCs" "Ptn" "Alk" "MuSCs" `
    3 1 3.975769333333333 11.927308 `
    0.03635351138648862 0.03635351138648861 `
    2 0.6666666666666666 0.012713 0.03813900000000001 `
    0.2207514079493428 0.2207514079493428 `
    0.05054395553466667 0.4548955998120001 `
    0.008025088822469826 0.008025088822469825

# Row 8 (new): Resolving-Mac -> Ptn -> Alk -> FAPs
Set-Row 8 "Resolving-Mac" "Ptn" "Alk" "FAPs" `
    1 0.3333333333333333 0.02199333333333333 0.06598 `
    0.0002011019319095741 0.0002011019319095741 `
    1 0.3333333333333333 0.04487666666666667 0.13463 `
    0.7792485920506572 0.7792485920506572 `
    0.0009869874888888889 0.008882887399999999 `
    0.0001567083972992028 0.0001567083972992028

# Row 9 (new): Resolving-Mac -> Ptn -> Alk -> MuSCs
Set-Row 9 "Resolving-Mac" "Ptn" "Alk" "MuSCs" `
    1 0.3333333333333333 0.02199333333333333 0.06598 `
    0.0002011019319095741 0.0002011019319095741 `
    2 0.6666666666666666 0.012713 0.03813900000000001 `
    0.2207514079493428 0.2207514079493428 `
    0.0002796012466666667 0.002516411220000001 `
    0.00004439353461037135 0.00004439353461037135
